$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column E ("reviews_count"), shifting columns F:K left to E:J
$ws.Range("E1").EntireColumn.Delete()
